# Horarios 141 - scrape refresh (16:36:34 -> 16:43:14)
# Updates rows that were re-ordered/re-scraped and appends newly scraped rows
# across the three worksheets (LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 (sheet1) ----
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 16:43:14"
$ws1.Range("A3").Value = "Total filas: 371"

# Cell-level updates for re-ordered / re-scraped rows
$ws1.Cells.Item(37, 1).Value = "06:33:46"
$ws1.Cells.Item(37, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(37, 4).Value = 63
$ws1.Cells.Item(38, 1).Value = "05:42:22"
$ws1.Cells.Item(38, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(38, 4).Value = 114
$ws1.Cells.Item(50, 1).Value = "07:12:53"
$ws1.Cells.Item(50, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(50, 4).Value = 49
$ws1.Cells.Item(51, 1).Value = "06:45:50"
$ws1.Cells.Item(51, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(51, 4).Value = 76
$ws1.Cells.Item(84, 1).Value = "07:36:59"
$ws1.Cells.Item(84, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(84, 4).Value = 88
$ws1.Cells.Item(85, 1).Value = "08:39:08"
$ws1.Cells.Item(85, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(85, 4).Value = 25
$ws1.Cells.Item(109, 1).Value = "08:46:25"
$ws1.Cells.Item(109, 3).Value = "14_ABASTO"
$ws1.Cells.Item(109, 4).Value = 78
$ws1.Cells.Item(110, 1).Value = "10:04:17"
$ws1.Cells.Item(110, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(110, 4).Value = 0
$ws1.Cells.Item(154, 1).Value = "10:04:17"
$ws1.Cells.Item(154, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(154, 4).Value = 77
$ws1.Cells.Item(155, 1).Value = "10:48:14"
$ws1.Cells.Item(155, 3).Value = "10_OLMOS"
$ws1.Cells.Item(155, 4).Value = 33
$ws1.Cells.Item(164, 1).Value = "11:34:25"
$ws1.Cells.Item(164, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(164, 4).Value = 0
$ws1.Cells.Item(165, 1).Value = "10:04:17"
$ws1.Cells.Item(165, 3).Value = "10_OLMOS"
$ws1.Cells.Item(165, 4).Value = 90
$ws1.Cells.Item(213, 1).Value = "11:11:31"
$ws1.Cells.Item(213, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(213, 4).Value = 103
$ws1.Cells.Item(214, 1).Value = "12:53:14"
$ws1.Cells.Item(214, 3).Value = "17_ROMERO"
$ws1.Cells.Item(214, 4).Value = 1
$ws1.Cells.Item(215, 1).Value = "12:32:47"
$ws1.Cells.Item(215, 3).Value = "10_OLMOS"
$ws1.Cells.Item(215, 4).Value = 22
$ws1.Cells.Item(224, 3).Value = "215_ALUAR"
$ws1.Cells.Item(225, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(236, 3).Value = "10_OLMOS"
$ws1.Cells.Item(237, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(258, 1).Value = "12:11:45"
$ws1.Cells.Item(258, 3).Value = "17_ROMERO"
$ws1.Cells.Item(258, 4).Value = 110
$ws1.Cells.Item(259, 1).Value = "13:51:48"
$ws1.Cells.Item(259, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(259, 4).Value = 10
$ws1.Cells.Item(283, 1).Value = "12:45:57"
$ws1.Cells.Item(283, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(283, 4).Value = 119
$ws1.Cells.Item(284, 1).Value = "14:44:53"
$ws1.Cells.Item(284, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(284, 4).Value = 0
$ws1.Cells.Item(293, 1).Value = "14:31:57"
$ws1.Cells.Item(293, 3).Value = "10_OLMOS"
$ws1.Cells.Item(293, 4).Value = 42
$ws1.Cells.Item(294, 1).Value = "13:39:24"
$ws1.Cells.Item(294, 3).Value = "14_ABASTO"
$ws1.Cells.Item(294, 4).Value = 94
$ws1.Cells.Item(329, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(330, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(345, 1).Value = "16:43:14"
$ws1.Cells.Item(345, 2).Value = "17:03"
$ws1.Cells.Item(345, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(345, 4).Value = 20
$ws1.Cells.Item(346, 1).Value = "15:17:21"
$ws1.Cells.Item(346, 2).Value = "17:04"
$ws1.Cells.Item(346, 3).Value = "14_ABASTO"
$ws1.Cells.Item(346, 4).Value = 107
$ws1.Cells.Item(347, 1).Value = "16:43:14"
$ws1.Cells.Item(347, 2).Value = "17:05"
$ws1.Cells.Item(347, 3).Value = "14_ABASTO"
$ws1.Cells.Item(347, 4).Value = 22
$ws1.Cells.Item(348, 2).Value = "17:07"
$ws1.Cells.Item(348, 3).Value = "15_ABASTO"
$ws1.Cells.Item(348, 4).Value = 82
$ws1.Cells.Item(349, 1).Value = "16:27:37"
$ws1.Cells.Item(349, 2).Value = "17:14"
$ws1.Cells.Item(349, 3).Value = "10_OLMOS"
$ws1.Cells.Item(349, 4).Value = 47
$ws1.Cells.Item(350, 1).Value = "15:45:31"
$ws1.Cells.Item(350, 2).Value = "17:17"
$ws1.Cells.Item(350, 3).Value = "17_ROMERO"
$ws1.Cells.Item(350, 4).Value = 92
$ws1.Cells.Item(351, 1).Value = "15:45:31"
$ws1.Cells.Item(351, 2).Value = "17:24"
$ws1.Cells.Item(351, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(351, 4).Value = 99
$ws1.Cells.Item(352, 1).Value = "16:43:14"
$ws1.Cells.Item(352, 2).Value = "17:25"
$ws1.Cells.Item(352, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(352, 4).Value = 42
$ws1.Cells.Item(353, 1).Value = "15:57:19"
$ws1.Cells.Item(353, 2).Value = "17:27"
$ws1.Cells.Item(353, 3).Value = "15_ABASTO"
$ws1.Cells.Item(353, 4).Value = 90
$ws1.Cells.Item(354, 1).Value = "16:36:34"
$ws1.Cells.Item(354, 2).Value = "17:30"
$ws1.Cells.Item(354, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(354, 4).Value = 54
$ws1.Cells.Item(355, 1).Value = "16:43:14"
$ws1.Cells.Item(355, 2).Value = "17:31"
$ws1.Cells.Item(355, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(355, 4).Value = 78
$ws1.Cells.Item(356, 2).Value = "17:34"
$ws1.Cells.Item(356, 3).Value = "10_OLMOS"
$ws1.Cells.Item(356, 4).Value = 109
$ws1.Cells.Item(357, 1).Value = "16:13:19"
$ws1.Cells.Item(357, 2).Value = "17:35"
$ws1.Cells.Item(357, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(357, 4).Value = 82
$ws1.Cells.Item(358, 1).Value = "15:45:31"
$ws1.Cells.Item(358, 2).Value = "17:35"
$ws1.Cells.Item(358, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(358, 4).Value = 110
$ws1.Cells.Item(359, 1).Value = "15:45:31"
$ws1.Cells.Item(359, 2).Value = "17:36"
$ws1.Cells.Item(359, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(359, 4).Value = 111
$ws1.Cells.Item(360, 1).Value = "16:43:14"
$ws1.Cells.Item(360, 2).Value = "17:37"
$ws1.Cells.Item(360, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(360, 4).Value = 54
$ws1.Cells.Item(361, 1).Value = "15:45:31"
$ws1.Cells.Item(361, 2).Value = "17:38"
$ws1.Cells.Item(361, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(361, 4).Value = 113
$ws1.Cells.Item(362, 1).Value = "16:36:34"
$ws1.Cells.Item(362, 2).Value = "17:44"
$ws1.Cells.Item(362, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(362, 4).Value = 77
$ws1.Cells.Item(363, 1).Value = "15:57:19"
$ws1.Cells.Item(363, 2).Value = "17:45"
$ws1.Cells.Item(363, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(363, 4).Value = 108
$ws1.Cells.Item(364, 1).Value = "16:43:14"
$ws1.Cells.Item(364, 2).Value = "17:47"
$ws1.Cells.Item(364, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(364, 4).Value = 64
$ws1.Cells.Item(365, 1).Value = "16:27:37"
$ws1.Cells.Item(365, 2).Value = "17:49"
$ws1.Cells.Item(365, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(365, 4).Value = 82
$ws1.Cells.Item(366, 1).Value = "15:57:19"
$ws1.Cells.Item(366, 2).Value = "17:51"
$ws1.Cells.Item(366, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(366, 4).Value = 114
$ws1.Cells.Item(367, 1).Value = "16:27:37"
$ws1.Cells.Item(367, 2).Value = "18:02"
$ws1.Cells.Item(367, 3).Value = "17_ROMERO"
$ws1.Cells.Item(367, 4).Value = 95

# New rows appended (368-376)
$ws1.Cells.Item(368, 1).Value = "16:13:19"
$ws1.Cells.Item(368, 2).Value = "18:03"
$ws1.Cells.Item(368, 3).Value = "17_ROMERO"
$ws1.Cells.Item(368, 4).Value = 110
$ws1.Cells.Item(368, 5).Value = "LP1912"
$ws1.Cells.Item(369, 1).Value = "16:13:19"
$ws1.Cells.Item(369, 2).Value = "18:04"
$ws1.Cells.Item(369, 3).Value = "14_ABASTO"
$ws1.Cells.Item(369, 4).Value = 111
$ws1.Cells.Item(369, 5).Value = "LP1912"
$ws1.Cells.Item(370, 1).Value = "16:36:34"
$ws1.Cells.Item(370, 2).Value = "18:04"
$ws1.Cells.Item(370, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(370, 4).Value = 88
$ws1.Cells.Item(370, 5).Value = "LP1912"
$ws1.Cells.Item(371, 1).Value = "16:43:14"
$ws1.Cells.Item(371, 2).Value = "18:05"
$ws1.Cells.Item(371, 3).Value = "14_ABASTO"
$ws1.Cells.Item(371, 4).Value = 82
$ws1.Cells.Item(371, 5).Value = "LP1912"
$ws1.Cells.Item(372, 1).Value = "16:27:37"
$ws1.Cells.Item(372, 2).Value = "18:24"
$ws1.Cells.Item(372, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(372, 4).Value = 117
$ws1.Cells.Item(372, 5).Value = "LP1912"
$ws1.Cells.Item(373, 1).Value = "16:43:14"
$ws1.Cells.Item(373, 2).Value = "18:25"
$ws1.Cells.Item(373, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(373, 4).Value = 102
$ws1.Cells.Item(373, 5).Value = "LP1912"
$ws1.Cells.Item(374, 1).Value = "16:36:34"
$ws1.Cells.Item(374, 2).Value = "18:34"
$ws1.Cells.Item(374, 3).Value = "14X44_ABASTO"
$ws1.Cells.Item(374, 4).Value = 118
$ws1.Cells.Item(374, 5).Value = "LP1912"
$ws1.Cells.Item(375, 1).Value = "16:43:14"
$ws1.Cells.Item(375, 2).Value = "18:38"
$ws1.Cells.Item(375, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(375, 4).Value = 115
$ws1.Cells.Item(375, 5).Value = "LP1912"
$ws1.Cells.Item(376, 1).Value = "16:43:14"
$ws1.Cells.Item(376, 2).Value = "18:41"
$ws1.Cells.Item(376, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(376, 4).Value = 118
$ws1.Cells.Item(376, 5).Value = "LP1912"

# ---- Sheet: LP1912-215 (sheet2) ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 16:43:14"

# ---- Sheet: 6203-6173 (sheet3) ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 16:43:14"
$ws3.Range("A3").Value = "Total filas: 43"

# New row 48
$ws3.Cells.Item(48, 1).Value = "16:43:14"
$ws3.Cells.Item(48, 2).Value = "18:22"
$ws3.Cells.Item(48, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(48, 4).Value = 99
$ws3.Cells.Item(48, 5).Value = "L6203"